$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 53: overall average of the three "211 trained"/res column averages ---
$ws.Range("L53").Formula = "=AVERAGE(L21,L33,L45,)"

# --- Row 56: header labels for the new "1RC"/"2RC"/"3RC" block ---
$ws.Cells.Item(56, 2).Value2 = "1RC"
$ws.Cells.Item(56, 3).Value2 = "2RC"
$ws.Cells.Item(56, 4).Value2 = "3RC"

# --- Rows 57-66: new data block (A = percentage level, B:D = values) ---
$data = @(
    @(100, 0.47239999999999999, 0.47003719999999999, 0.29702319999999999),
    @(90,  0.38482,             0.45051999999999998, 0.25930999999999998),
    @(80,  0.39062999999999998, 0.48894700000000002, 0.2723777),
    @(70,  0.40162999999999999, 0.490838,             0.29239599999999999),
    @(60,  0.40132200000000001, 0.48226799999999997, 0.30452000000000001),
    @(50,  0.38605099999999998, 0.48337000000000002, 0.28462799999999999),
    @(40,  0.41930800000000001, 0.57841480000000001, 0.36316599999999999),
    @(30,  0.4473008,           0.66423299999999996, 0.41869469999999998),
    @(20,  0.43763000000000002, 0.57943960000000005, 0.40178950000000002),
    @(10,  0.38205040000000001, 0.21379100000000001, 0.21379000000000001)
)

$row = 57
foreach ($vals in $data) {
    $ws.Cells.Item($row, 1).Value2 = $vals[0]
    $ws.Cells.Item($row, 2).Value2 = $vals[1]
    $ws.Cells.Item($row, 3).Value2 = $vals[2]
    $ws.Cells.Item($row, 4).Value2 = $vals[3]
    $row = $row + 1
}

# First data row (57) carries the scientific-notation number format on B:D.
$ws.Range("B57:D57").NumberFormat = "0.00E+00"

# --- Row 68: column averages for the new block ---
$ws.Range("B68").Formula = "=AVERAGE(B57:B66)"
$ws.Range("D68").Formula = "=AVERAGE(D57:D66)"
$ws.Range("B68").NumberFormat = "0.00E+00"
$ws.Range("D68").NumberFormat = "0.00E+00"

# --- View state: selection moves to D69 as in the saved workbook ---
$ws.Range("D69").Select()
